$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.055275000000001
$ws.Range("H2").Value = 12.165825
$ws.Range("I2").Value = 0.1947228515851206
$ws.Range("J2").Value = 0.1947228515851206
$ws.Range("Q2").Value = 1.113689359183334
$ws.Range("R2").Value = 10.02320423265
$ws.Range("S2").Value = 0.1947228515851206
$ws.Range("T2").Value = 0.1947228515851206

$ws.Range("G3").Value = 6.542812333333333
$ws.Range("I3").Value = 0.3141673684110111
$ws.Range("J3").Value = 0.3141673684110111
$ws.Range("S3").Value = 0.3141673684110111
$ws.Range("T3").Value = 0.3141673684110111

$ws.Range("G4").Value = 6.7998
$ws.Range("H4").Value = 20.3994
$ws.Range("I4").Value = 0.3265071903159472
$ws.Range("J4").Value = 0.3265071903159472
$ws.Range("Q4").Value = 1.8674109412
$ws.Range("R4").Value = 16.8066984708
$ws.Range("S4").Value = 0.3265071903159472
$ws.Range("T4").Value = 0.3265071903159472

$ws.Range("G5").Value = 3.427994
$ws.Range("H5").Value = 10.283982
$ws.Range("I5").Value = 0.164602589687921
$ws.Range("J5").Value = 0.164602589687921
$ws.Range("Q5").Value = 0.9414208509026667
$ws.Range("R5").Value = 8.472787658124
$ws.Range("S5").Value = 0.164602589687921
$ws.Range("T5").Value = 0.164602589687921
